$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "28.01.2019"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "Enemy fights back now yay. It's almost a game."

$ws.Range("C10").Select()
